$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").Formula = "'328.62"
$ws.Range("D2").Style = $s
$s = $ws.Range("E2").Style
$ws.Range("E2").Formula = "'1.40%"
$ws.Range("E2").Style = $s
$s = $ws.Range("D3").Style
$ws.Range("D3").Formula = "'41.65"
$ws.Range("D3").Style = $s
$s = $ws.Range("E3").Style
$ws.Range("E3").Formula = "'5.22%"
$ws.Range("E3").Style = $s
$s = $ws.Range("D4").Style
$ws.Range("D4").Formula = "'5.626"
$ws.Range("D4").Style = $s
$s = $ws.Range("E4").Style
$ws.Range("E4").Formula = "'-4.04%"
$ws.Range("E4").Style = $s
$s = $ws.Range("E5").Style
$ws.Range("E5").Formula = "'1.79%"
$ws.Range("E5").Style = $s
$s = $ws.Range("D6").Style
$ws.Range("D6").Formula = "'2.025"
$ws.Range("D6").Style = $s
$s = $ws.Range("E6").Style
$ws.Range("E6").Formula = "'0.31%"
$ws.Range("E6").Style = $s
$s = $ws.Range("D7").Style
$ws.Range("D7").Formula = "'8.737"
$ws.Range("D7").Style = $s
$s = $ws.Range("E7").Style
$ws.Range("E7").Formula = "'1.26%"
$ws.Range("E7").Style = $s
$s = $ws.Range("D8").Style
$ws.Range("D8").Formula = "'4.524"
$ws.Range("D8").Style = $s
$s = $ws.Range("E8").Style
$ws.Range("E8").Formula = "'-1.01%"
$ws.Range("E8").Style = $s
$s = $ws.Range("E9").Style
$ws.Range("E9").Formula = "'0.06%"
$ws.Range("E9").Style = $s
$s = $ws.Range("D10").Style
$ws.Range("D10").Formula = "'0.9217"
$ws.Range("D10").Style = $s
$s = $ws.Range("E10").Style
$ws.Range("E10").Formula = "'-0.86%"
$ws.Range("E10").Style = $s
$s = $ws.Range("D11").Style
$ws.Range("D11").Formula = "'0.1275"
$ws.Range("D11").Style = $s
$s = $ws.Range("E11").Style
$ws.Range("E11").Formula = "'0.59%"
$ws.Range("E11").Style = $s
$s = $ws.Range("D12").Style
$ws.Range("D12").Formula = "'0.1956"
$ws.Range("D12").Style = $s
$s = $ws.Range("E12").Style
$ws.Range("E12").Formula = "'-0.08%"
$ws.Range("E12").Style = $s
$s = $ws.Range("D13").Style
$ws.Range("D13").Formula = "'0.09310"
$ws.Range("D13").Style = $s
$s = $ws.Range("E13").Style
$ws.Range("E13").Formula = "'1.80%"
$ws.Range("E13").Style = $s
$s = $ws.Range("D14").Style
$ws.Range("D14").Formula = "'0.03810"
$ws.Range("D14").Style = $s
$s = $ws.Range("E14").Style
$ws.Range("E14").Formula = "'6.24%"
$ws.Range("E14").Style = $s
$s = $ws.Range("E15").Style
$ws.Range("E15").Formula = "'0.87%"
$ws.Range("E15").Style = $s
$s = $ws.Range("D16").Style
$ws.Range("D16").Formula = "'0.001305"
$ws.Range("D16").Style = $s
$s = $ws.Range("E16").Style
$ws.Range("E16").Formula = "'0.93%"
$ws.Range("E16").Style = $s
$s = $ws.Range("D17").Style
$ws.Range("D17").Formula = "'0.006286"
$ws.Range("D17").Style = $s
$s = $ws.Range("E17").Style
$ws.Range("E17").Formula = "'2.09%"
$ws.Range("E17").Style = $s
$s = $ws.Range("D19").Style
$ws.Range("D19").Formula = "'3.438"
$ws.Range("D19").Style = $s
$s = $ws.Range("E19").Style
$ws.Range("E19").Formula = "'2.57%"
$ws.Range("E19").Style = $s
$s = $ws.Range("D21").Style
$ws.Range("D21").Formula = "'8.327"
$ws.Range("D21").Style = $s
$s = $ws.Range("E21").Style
$ws.Range("E21").Formula = "'-4.32%"
$ws.Range("E21").Style = $s
$s = $ws.Range("D22").Style
$ws.Range("D22").Formula = "'0.1394"
$ws.Range("D22").Style = $s
$s = $ws.Range("E22").Style
$ws.Range("E22").Formula = "'1.73%"
$ws.Range("E22").Style = $s
$s = $ws.Range("E23").Style
$ws.Range("E23").Formula = "'-1.40%"
$ws.Range("E23").Style = $s
$s = $ws.Range("E24").Style
$ws.Range("E24").Formula = "'0.03%"
$ws.Range("E24").Style = $s
$s = $ws.Range("D25").Style
$ws.Range("D25").Formula = "'0.001260"
$ws.Range("D25").Style = $s
$s = $ws.Range("E25").Style
$ws.Range("E25").Formula = "'-0.16%"
$ws.Range("E25").Style = $s
$s = $ws.Range("D26").Style
$ws.Range("D26").Formula = "'0.004338"
$ws.Range("D26").Style = $s
$s = $ws.Range("E26").Style
$ws.Range("E26").Formula = "'-1.24%"
$ws.Range("E26").Style = $s
$s = $ws.Range("D27").Style
$ws.Range("D27").Formula = "'0.0001181"
$ws.Range("D27").Style = $s
$s = $ws.Range("E27").Style
$ws.Range("E27").Formula = "'2.74%"
$ws.Range("E27").Style = $s
$s = $ws.Range("D39").Style
$ws.Range("D39").Formula = "'0.02772"
$ws.Range("D39").Style = $s
$s = $ws.Range("E39").Style
$ws.Range("E39").Formula = "'9.82%"
$ws.Range("E39").Style = $s
$s = $ws.Range("D40").Style
$ws.Range("D40").Formula = "'0.05422"
$ws.Range("D40").Style = $s
$s = $ws.Range("E40").Style
$ws.Range("E40").Formula = "'3.39%"
$ws.Range("E40").Style = $s
$s = $ws.Range("D41").Style
$ws.Range("D41").Formula = "'0.007669"
$ws.Range("D41").Style = $s
$s = $ws.Range("E41").Style
$ws.Range("E41").Formula = "'2.84%"
$ws.Range("E41").Style = $s
$s = $ws.Range("D42").Style
$ws.Range("D42").Formula = "'0.1418"
$ws.Range("D42").Style = $s
$s = $ws.Range("E42").Style
$ws.Range("E42").Formula = "'0.75%"
$ws.Range("E42").Style = $s
$s = $ws.Range("D43").Style
$ws.Range("D43").Formula = "'0.008981"
$ws.Range("D43").Style = $s
$s = $ws.Range("E43").Style
$ws.Range("E43").Formula = "'-6.48%"
$ws.Range("E43").Style = $s
$s = $ws.Range("D44").Style
$ws.Range("D44").Formula = "'0.002132"
$ws.Range("D44").Style = $s
$s = $ws.Range("E44").Style
$ws.Range("E44").Formula = "'0.77%"
$ws.Range("E44").Style = $s
$s = $ws.Range("D45").Style
$ws.Range("D45").Formula = "'0.01170"
$ws.Range("D45").Style = $s
$s = $ws.Range("E45").Style
$ws.Range("E45").Formula = "'17.19%"
$ws.Range("E45").Style = $s
$s = $ws.Range("D46").Style
$ws.Range("D46").Formula = "'0.00006649"
$ws.Range("D46").Style = $s
$s = $ws.Range("E46").Style
$ws.Range("E46").Formula = "'-1.35%"
$ws.Range("E46").Style = $s
$s = $ws.Range("D48").Style
$ws.Range("D48").Formula = "'0.003214"
$ws.Range("D48").Style = $s
$s = $ws.Range("E48").Style
$ws.Range("E48").Formula = "'7.13%"
$ws.Range("E48").Style = $s
$s = $ws.Range("D50").Style
$ws.Range("D50").Formula = "'0.00002104"
$ws.Range("D50").Style = $s

Write-Output "done"
